# New Microsite scripts support to Beta server
# Appends newly-recorded certificate-registration test runs to the
# "AMSIN" history sheet and the "AMS" history sheet, and fixes up a
# previously mis-entered timestamp / missing style on AMS row 17.

$wb = $excel.ActiveWorkbook

$amsin = $wb.Worksheets.Item("AMSIN")
$ams   = $wb.Worksheets.Item("AMS")

# ---------------------------------------------------------------------
# AMSIN: append rows 60-63 (dimension grows from A1:G59 to A1:G63)
# ---------------------------------------------------------------------
$amsinRows = @(
    @{ Row = 60; Date = "2022-09-15"; Time = 44819.62238364583; Name = "certiecs166";   Total = 51; Pass = 51; Fail = 0; Taken = 1.53 },
    @{ Row = 61; Date = "2022-09-16"; Time = 44820.64392523148; Name = "fstc167";       Total = 51; Pass = 51; Fail = 0; Taken = 1.07 },
    @{ Row = 62; Date = "2022-09-19"; Time = 44823.60025365741; Name = "scndcerti167";  Total = 51; Pass = 50; Fail = 1; Taken = 2.57 },
    @{ Row = 63; Date = "2022-09-20"; Time = 44824.38537011574; Name = "finalcerti167"; Total = 51; Pass = 51; Fail = 0; Taken = 1.46 }
)

foreach ($r in $amsinRows) {
    $rowNum = $r.Row

    # Column A holds a date-look-alike string ("2022-09-15"), not a real
    # date serial - a leading apostrophe forces literal text entry (same
    # General number format as the rest of the row) instead of Excel
    # auto-converting it to a date serial.
    $amsin.Cells.Item($rowNum, 1).Value = "'" + $r.Date

    # Column B is the real run timestamp (serial date/time number).
    $amsin.Cells.Item($rowNum, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $amsin.Cells.Item($rowNum, 2).Value = $r.Time

    $amsin.Cells.Item($rowNum, 3).Value = $r.Name
    $amsin.Cells.Item($rowNum, 4).Value = $r.Total
    $amsin.Cells.Item($rowNum, 5).Value = $r.Pass
    $amsin.Cells.Item($rowNum, 6).Value = $r.Fail
    $amsin.Cells.Item($rowNum, 7).Value = $r.Taken
}

# ---------------------------------------------------------------------
# AMS: correct the run-time precision on row 17, then append row 18
# (dimension grows from A1:G17 to A1:G18)
# ---------------------------------------------------------------------
$ams.Cells.Item(17, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ams.Cells.Item(17, 2).Value = 44812.50833604167

$ams.Cells.Item(18, 1).Value = "'2022-09-14"

$ams.Cells.Item(18, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ams.Cells.Item(18, 2).Value = 44818.98614082176

$ams.Cells.Item(18, 3).Value = "certhftx166"
$ams.Cells.Item(18, 4).Value = 51
$ams.Cells.Item(18, 5).Value = 51
$ams.Cells.Item(18, 6).Value = 0
$ams.Cells.Item(18, 7).Value = 0.99
